# Fruta / hortaliza, semanal
# Updates the weekly Breva (Vega Central Mapocho de Santiago) records:
# corrected dates, volumes, prices, quality, unit and origin fields per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Fecha 44553 -> 44558 ; Volumen 200 -> 20
$ws.Range("D2").Value = 44558
$ws.Range("M2").Value = 20

# Row 3: Fecha 44553 -> 44558 ; Volumen 150 -> 25
$ws.Range("D3").Value = 44558
$ws.Range("M3").Value = 25

# Row 4: Fecha 44558 -> 44561 ; Calidad Especial -> Primera ;
#        Volumen 20 -> 200 ; Precios 22000 -> 18000 ; Precio/Kg 3667 -> 3000
$ws.Range("D4").Value = 44561
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("S4").Value = 3000

# Row 5: Fecha 44558 -> 44550 ; Volumen 25 -> 60 ; Precios 18000 -> 24000 ;
#        Unidad "$/bandeja 6 kilos" -> "$/bandeja 7 kilos" ;
#        Origen Provincia de San Felipe de Aconcagua -> Región Metropolitana ;
#        Precio/Kg 3000 -> 3429 ; Kg/unidad 6 -> 7
$ws.Range("D5").Value = 44550
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 24000
$ws.Range("O5").Value = 24000
$ws.Range("P5").Value = 24000
$ws.Range("Q5").Value = "$/bandeja 7 kilos"
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 3429
$ws.Range("T5").Value = 7

# Row 7: Fecha 44187 -> 44553 ; Volumen 45 -> 200 ; Precios 14000 -> 22000 ;
#        Unidad "$/bandeja 7 kilos" -> "$/bandeja 6 kilos" ;
#        Precio/Kg 2000 -> 3667 ; Kg/unidad 7 -> 6
$ws.Range("D7").Value = 44553
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 22000
$ws.Range("O7").Value = 22000
$ws.Range("P7").Value = 22000
$ws.Range("Q7").Value = "$/bandeja 6 kilos"
$ws.Range("S7").Value = 3667
$ws.Range("T7").Value = 6

# Row 8: Fecha 44187 -> 44553 ; Volumen 50 -> 150 ; Precios 12000 -> 18000 ;
#        Unidad "$/bandeja 7 kilos" -> "$/bandeja 6 kilos" ;
#        Precio/Kg 1714 -> 3000 ; Kg/unidad 7 -> 6
$ws.Range("D8").Value = 44553
$ws.Range("M8").Value = 150
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 18000
$ws.Range("Q8").Value = "$/bandeja 6 kilos"
$ws.Range("S8").Value = 3000
$ws.Range("T8").Value = 6

# Row 9: Fecha 44561 -> 44572 ; Volumen 200 -> 65 ; Precios 18000 -> 20000 ;
#        Origen Provincia de San Felipe de Aconcagua -> Región Metropolitana ;
#        Precio/Kg 3000 -> 3333
$ws.Range("D9").Value = 44572
$ws.Range("M9").Value = 65
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("R9").Value = "Región Metropolitana"
$ws.Range("S9").Value = 3333

# Row 12: Fecha 44550 -> 44187 ; Calidad Primera -> Especial ;
#         Volumen 60 -> 45 ; Precios 24000 -> 14000 ;
#         Origen Región Metropolitana -> Provincia de San Felipe de Aconcagua ;
#         Precio/Kg 3429 -> 2000
$ws.Range("D12").Value = 44187
$ws.Range("L12").Value = "Especial"
$ws.Range("M12").Value = 45
$ws.Range("N12").Value = 14000
$ws.Range("O12").Value = 14000
$ws.Range("P12").Value = 14000
$ws.Range("R12").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S12").Value = 2000

# Row 13: Fecha 44572 -> 44187 ; Volumen 65 -> 50 ; Precios 20000 -> 12000 ;
#         Unidad "$/bandeja 6 kilos" -> "$/bandeja 7 kilos" ;
#         Origen Región Metropolitana -> Provincia de San Felipe de Aconcagua ;
#         Precio/Kg 3333 -> 1714 ; Kg/unidad 6 -> 7
$ws.Range("D13").Value = 44187
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 12000
$ws.Range("Q13").Value = "$/bandeja 7 kilos"
$ws.Range("R13").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S13").Value = 1714
$ws.Range("T13").Value = 7
